$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet lists daily price observations (Apio, Feria Lagunitas de Puerto
# Montt) in rows 2..158. A new observation is inserted at row 126, which
# pushes every following observation (127..158) down by one row, and the
# former last row (158) becomes a brand new row 159.
#
# Strategy:
#  1. Duplicate row 158 entirely into the brand new row 159 (before any
#     other change, so it captures the original, pre-edit data).
#  2. Walk rows 158 down to 127 (descending order) and copy the variable
#     columns (D, I, J, K, L, M, P) from the row immediately above. Working
#     top-down-in-reverse (i.e. from the bottom up) guarantees that the
#     "source" row for each copy has not yet been overwritten.
#  3. Finally, overwrite row 126 with the brand new observation's Date (D)
#     and Volumen (J) values; its other columns keep their original values.

# Step 1: clone row 158 -> row 159 (full row, preserves styles/number formats)
$ws.Range("A158:R158").Copy($ws.Range("A159:R159"))

# Step 2: shift rows 127..158 down from their predecessor, from bottom to top
for ($r = 158; $r -ge 127; $r--) {
    $prev = $r - 1
    $ws.Range("D$r").Value = $ws.Range("D$prev").Value2
    $ws.Range("I$r").Value = $ws.Range("I$prev").Value2
    $ws.Range("J$r").Value = $ws.Range("J$prev").Value2
    $ws.Range("K$r").Value = $ws.Range("K$prev").Value2
    $ws.Range("L$r").Value = $ws.Range("L$prev").Value2
    $ws.Range("M$r").Value = $ws.Range("M$prev").Value2
    $ws.Range("P$r").Value = $ws.Range("P$prev").Value2
}

# Step 3: set the new observation inserted at row 126 (Date + Volumen only)
$ws.Range("D126").Value = 44511
$ws.Range("J126").Value = 25
